$wb = $excel.ActiveWorkbook

# --- Update "Latest Handoff Datetime" / generation timestamps ---
# Overview sheet: "Latest HO Xliff Generate Date" column G, rows 8-13
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G8:G13").Value = "2016-08-13 04:24:51"

# zh-cn sheet: "Latest Handoff Datetime" column H, rows 8-13
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H8:H13").Value = "2016-08-13 04:24:44"

# de-de sheet: "Latest Handoff Datetime" column H, rows 8-13
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H8:H13").Value = "2016-08-13 04:24:51"

# --- Update "Priority" column E, rows 8-13, to "ht" (handoff type) ---
# These rows are the "Ready for handoff" rows being handed off in this report.
$wsZhCn.Range("E8:E13").Value = "ht"
$wsDeDe.Range("E8:E13").Value = "ht"
